$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1607.3448
$ws.Range("J17").Value = 1607.3448
$ws.Range("L17").Value = 4822.0344
$ws.Range("N17").Value = -5158.0344
$ws.Range("H70").Value = 1617
$ws.Range("J70").Value = 1646.25
$ws.Range("L70").Value = 4938.75
$ws.Range("N70").Value = -5478.75
$ws.Range("H73").Value = 1617
$ws.Range("J73").Value = 1646.25
$ws.Range("L73").Value = 4938.75
$ws.Range("N73").Value = -6810.75
$ws.Range("H100").Value = 2199.875
$ws.Range("I100").Value = 2199.875
$ws.Range("K100").Value = 2199.875
$ws.Range("M100").Value = -1658.875
$ws.Range("H112").Value = 4242
$ws.Range("I112").Value = 2000
$ws.Range("J112").Value = 4366.5557
$ws.Range("K112").Value = 6000
$ws.Range("L112").Value = 13099.6671
$ws.Range("M112").Value = -4892
$ws.Range("N112").Value = -15315.6671
$ws.Range("H129").Value = 858.5323
$ws.Range("J129").Value = 856.6491
$ws.Range("L129").Value = 2569.9473
$ws.Range("N129").Value = -12569.9473
$ws.Range("H135").Value = 682.2941
$ws.Range("I135").Value = 584.9231
$ws.Range("K135").Value = 5264.3079
$ws.Range("M135").Value = -2729.3079
$ws.Range("H140").Value = 81835.31
$ws.Range("J140").Value = 81835.31
$ws.Range("L140").Value = 81835.31
$ws.Range("N140").Value = -92195.31

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1158.875
$ws.Range("I2").Value = 1204.2
$ws.Range("J2").Value = 1083.3334
$ws.Range("K2").Value = 1204.2
$ws.Range("L2").Value = 1083.3334
$ws.Range("M2").Value = -1091.2
$ws.Range("N2").Value = -1309.3334
$ws.Range("H32").Value = 2226.602
$ws.Range("I32").Value = 1489.9048
$ws.Range("J32").Value = 6646.7856
$ws.Range("K32").Value = 1489.9048
$ws.Range("L32").Value = 6646.7856
$ws.Range("M32").Value = -1202.9048
$ws.Range("N32").Value = -7220.7856
$ws.Range("H74").Value = 1569.3125
$ws.Range("I74").Value = 1408.909
$ws.Range("K74").Value = 1408.909
$ws.Range("M74").Value = -534.9090000000001
$ws.Range("H77").Value = 1569.3125
$ws.Range("I77").Value = 1408.909
$ws.Range("K77").Value = 7044.545
$ws.Range("M77").Value = -2676.545
$ws.Range("H116").Value = 1158.875
$ws.Range("I116").Value = 1204.2
$ws.Range("J116").Value = 1083.3334
$ws.Range("K116").Value = 1204.2
$ws.Range("L116").Value = 1083.3334
$ws.Range("M116").Value = 1089.8
$ws.Range("N116").Value = -5671.3334

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1158.875
$ws.Range("I3").Value = 1204.2
$ws.Range("J3").Value = 1083.3334
$ws.Range("K3").Value = 1204.2
$ws.Range("L3").Value = 1083.3334
$ws.Range("M3").Value = -1090.2
$ws.Range("N3").Value = -1311.3334
$ws.Range("H22").Value = 190
$ws.Range("I22").Value = 190
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 190
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -17
$ws.Range("N22").ClearContents()
$ws.Range("H94").Value = 1047.1818
$ws.Range("I94").Value = 835.3333
$ws.Range("K94").Value = 835.3333
$ws.Range("M94").Value = -384.3333
$ws.Range("H105").Value = 2513.2
$ws.Range("J105").Value = 2666.3333
$ws.Range("L105").Value = 2666.3333
$ws.Range("N105").Value = -6160.3333
$ws.Range("H107").Value = 7500.5
$ws.Range("I107").Value = 7500.5
$ws.Range("K107").Value = 7500.5
$ws.Range("M107").Value = -5580.5
$ws.Range("H134").Value = 6333.615
$ws.Range("I134").Value = 7638.7
$ws.Range("K134").Value = 22916.1
$ws.Range("M134").Value = -20381.1

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1856.7368
$ws.Range("J31").Value = 2666
$ws.Range("L31").Value = 2666
$ws.Range("N31").Value = -3256
$ws.Range("H34").Value = 1856.7368
$ws.Range("J34").Value = 2666
$ws.Range("L34").Value = 2666
$ws.Range("N34").Value = -3070
$ws.Range("H58").Value = 1662.6897
$ws.Range("I58").Value = 1009.5238
$ws.Range("K58").Value = 1009.5238
$ws.Range("M58").Value = -806.5238000000001
$ws.Range("H132").Value = 3067.2222
$ws.Range("J132").Value = 4141.2
$ws.Range("L132").Value = 12423.6
$ws.Range("N132").Value = -17483.6
$ws.Range("H134").Value = 2378.4736
$ws.Range("I134").Value = 1955.1111
$ws.Range("J134").Value = 9999
$ws.Range("K134").Value = 5865.3333
$ws.Range("L134").Value = 29997
$ws.Range("M134").Value = -3330.3333
$ws.Range("N134").Value = -35067
$ws.Range("H136").Value = 1662.6897
$ws.Range("I136").Value = 1009.5238
$ws.Range("K136").Value = 3028.5714
$ws.Range("M136").Value = -478.5714000000003

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 39.333332
$ws.Range("I8").Value = 39.333332
$ws.Range("K8").Value = 117.999996
$ws.Range("M8").Value = 21.000004
$ws.Range("H56").Value = 10008.125
$ws.Range("I56").Value = 10008.125
$ws.Range("K56").Value = 10008.125
$ws.Range("M56").Value = -9478.125
$ws.Range("H60").Value = 3000
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H131").Value = 1237.15
$ws.Range("I131").Value = 609.5
$ws.Range("J131").Value = 1277.2128
$ws.Range("K131").Value = 1828.5
$ws.Range("L131").Value = 3831.6384
$ws.Range("M131").Value = 3211.5
$ws.Range("N131").Value = -13911.6384

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1640.6
$ws.Range("I107").Value = 300
$ws.Range("J107").Value = 2534.3333
$ws.Range("K107").Value = 300
$ws.Range("L107").Value = 2534.3333
$ws.Range("M107").Value = 1620
$ws.Range("N107").Value = -6374.3333
$ws.Range("H122").Value = 1824.174
$ws.Range("J122").Value = 2466.3333
$ws.Range("L122").Value = 7398.999899999999
$ws.Range("N122").Value = -12298.9999
$ws.Range("H126").Value = 46564.87
$ws.Range("I126").Value = 3226.9092
$ws.Range("K126").Value = 9680.7276
$ws.Range("M126").Value = -7210.7276

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3995
$ws.Range("J82").Value = 4660
$ws.Range("L82").Value = 4660
$ws.Range("N82").Value = -5382
$ws.Range("H85").Value = 3995
$ws.Range("J85").Value = 4660
$ws.Range("L85").Value = 4660
$ws.Range("N85").Value = -7156
$ws.Range("H132").Value = 3969.4
$ws.Range("I132").Value = 2633
$ws.Range("K132").Value = 7899
$ws.Range("M132").Value = -5369
$ws.Range("H136").Value = 3756
$ws.Range("I136").Value = 2876.8
$ws.Range("J136").Value = 4855
$ws.Range("K136").Value = 8630.400000000001
$ws.Range("L136").Value = 14565
$ws.Range("M136").Value = -6080.400000000001
$ws.Range("N136").Value = -19665

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 14000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 14000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 14000
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -14346
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H81").Value = 2079.3333
$ws.Range("I81").Value = 2079.3333
$ws.Range("K81").Value = 4158.6666
$ws.Range("M81").Value = -3097.6666
$ws.Range("H84").Value = 2079.3333
$ws.Range("I84").Value = 2079.3333
$ws.Range("K84").Value = 20793.333
$ws.Range("M84").Value = -15489.333
